# Fruta / hortaliza, semanal
# Insert two new weekly price rows (date 44516, Mandarina - Murcott,
# Primera/Segunda) at the top of the Feria Lagunitas de Puerto Montt
# table, pushing the existing rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 122.
$ws.Range("A122:A123").EntireRow.Insert()

# New row 122: Mandarina / Murcott / Primera
$ws.Range("A122").Value = 4
$ws.Range("B122").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C122").Value = "Los Lagos"
$ws.Range("D122").Value = 44516
$ws.Range("E122").Value = 10
$ws.Range("F122").Value = "Fruta"
$ws.Range("G122").Value = 100102
$ws.Range("H122").Value = "Cítricos"
$ws.Range("I122").Value = 100102004
$ws.Range("J122").Value = "Mandarina"
$ws.Range("K122").Value = "Murcott"
$ws.Range("L122").Value = "Primera"
$ws.Range("M122").Value = 600
$ws.Range("N122").Value = 6500
$ws.Range("O122").Value = 7000
$ws.Range("P122").Value = 6750
$ws.Range("Q122").Value = "$/bandeja 10 kilos"
$ws.Range("R122").Value = "Provincia de Limarí"
$ws.Range("S122").Value = 675
$ws.Range("T122").Value = 10

# New row 123: Mandarina / Murcott / Segunda
$ws.Range("A123").Value = 4
$ws.Range("B123").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C123").Value = "Los Lagos"
$ws.Range("D123").Value = 44516
$ws.Range("E123").Value = 10
$ws.Range("F123").Value = "Fruta"
$ws.Range("G123").Value = 100102
$ws.Range("H123").Value = "Cítricos"
$ws.Range("I123").Value = 100102004
$ws.Range("J123").Value = "Mandarina"
$ws.Range("K123").Value = "Murcott"
$ws.Range("L123").Value = "Segunda"
$ws.Range("M123").Value = 300
$ws.Range("N123").Value = 5000
$ws.Range("O123").Value = 5000
$ws.Range("P123").Value = 5000
$ws.Range("Q123").Value = "$/bandeja 10 kilos"
$ws.Range("R123").Value = "Provincia de Limarí"
$ws.Range("S123").Value = 500
$ws.Range("T123").Value = 10
